$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows at row 14, pushing the old rows 14-17 ("D4","D7","D4T","D7T" tensile
# block) down to rows 18-21. Seed the new rows with the same formatting pattern used by
# the existing creep block in rows 10-13 (styles 11/12,13/14,13/14,15/16) by copying it.
$ws.Rows("14:17").Insert()
$ws.Range("A10:G13").Copy($ws.Range("A14:G17"))

# Fill in the new material rows' actual data.
$ws.Range("A14").Value = "G39"
$ws.Range("B14").Value = 617
$ws.Range("C14").Value = 1000
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = "creep"
$ws.Range("F14").Value = "time"
$ws.Range("G14").Value = "strain"

$ws.Range("A15").Value = "G52"
$ws.Range("B15").Value = 617
$ws.Range("C15").Value = 1000
$ws.Range("D15").Value = 12
$ws.Range("E15").Value = "creep"
$ws.Range("F15").Value = "time"
$ws.Range("G15").Value = "strain"

$ws.Range("A16").Value = "G30"
$ws.Range("B16").Value = 617
$ws.Range("C16").Value = 1000
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = "creep"
$ws.Range("F16").Value = "time"
$ws.Range("G16").Value = "strain"

$ws.Range("A17").Value = "G18"
$ws.Range("B17").Value = 617
$ws.Range("C17").Value = 1000
$ws.Range("D17").Value = 16
$ws.Range("E17").Value = "creep"
$ws.Range("F17").Value = "time"
$ws.Range("G17").Value = "strain"

# Fix the pre-existing B11/B12 style glitch (they used an orphan "fill3+border0" style
# rather than the plain "fill3" style used by the rest of that row's block); Excel's
# repair round trip folds them onto the same style used by C11:G11 / C12:G12.
$ws.Range("B11").Style = $ws.Range("C11").Style
$ws.Range("B12").Style = $ws.Range("C12").Style

# Update window/view metadata to match the refreshed sheet.
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("G14").Select()
